$d = $word.ActiveDocument

$old = "môžete pozorovať súhvezdie Ozvezdje Dvojčka 2022: 14.-23. februar, 14.-24. marec"
$new = "2022: Datumi kampanje za opazovanje Ozvezdje Dvojčka: 14.-23. februar, 14.-24. marec"

$range = $d.Content
$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
